# Updated viz for day 22.
# - Rename shared string "Day 22: TITLE" to the real title "Day 22: Slam Shuffle"
#   (done by writing the new text into the cell that held the placeholder; the
#   engine will append the new string and drop the now-unused placeholder,
#   renumbering every other reference the same way the source diff shows).
# - Fill in that day's results (row 26) and the Day 21 row (row 25) that had
#   been missing its "Finish"/rank data.
# - Move the sheet selection to H27.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2019")

# --- Row 25: Day 21 (Springdroid Adventure) results ---
$ws.Range("C25").Value = 0.026770833333333331
$ws.Range("E25").Value = 0.036724537037037035
$ws.Range("H25").Value = "4th"

# --- Row 26: Day 22 (Slam Shuffle) - give it its real title and results ---
$ws.Range("B26").Value = "Day 22: Slam Shuffle"
$ws.Range("C26").Value = 0.022685185185185183
$ws.Range("E26").Value = 0.14959490740740741
$ws.Range("F26").Value = 0.1292939814814815
$ws.Range("H26").Value = "2nd"

# --- Move active selection ---
$ws.Range("H27").Select()
